# Insert a new weekly price record as row 27 in the "Ajo" (garlic) sheet.
# This pushes the previous rows 27-30 down to rows 28-31 (data unchanged),
# and fills the freed-up row 27 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 27; existing rows 27-30 shift to 28-31.
$ws.Rows.Item(27).Insert()

# Populate the new row 27 with the new weekly record.
$ws.Range("A27").Value = 1
$ws.Range("B27").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C27").Value = "Arica y Parinacota"
$ws.Range("D27").Value = 44918
$ws.Range("D27").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E27").Value = 15
$ws.Range("F27").Value = 100112003
$ws.Range("G27").Value = "Ajo"
$ws.Range("H27").Value = "Chino"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 200
$ws.Range("K27").Value = 12000
$ws.Range("L27").Value = 13000
$ws.Range("M27").Value = 12250
$ws.Range("N27").Value = "`$/caja 10 kilos"
$ws.Range("O27").Value = "China"
$ws.Range("P27").Value = 1225
$ws.Range("Q27").Value = 10
$ws.Range("R27").Value = "Hortaliza"

Write-Host "Inserted new row 27; dimension now:" $ws.UsedRange.Address()
